# Update the "Time" column (J) on the "test-result" sheet for rows 2-7,
# replacing the old run timestamp with the new one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test-result")

$newTime = "2025-06-07 13:13:38"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 10).Value = $newTime
}
